$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1911.4286
$ws.Range("J40").Value = 3300
$ws.Range("L40").Value = 3300
$ws.Range("N40").Value = -3650
$ws.Range("H116").Value = 3596.4138
$ws.Range("I116").Value = 2720.2222
$ws.Range("J116").Value = 5030.1816
$ws.Range("K116").Value = 2720.2222
$ws.Range("L116").Value = 5030.1816
$ws.Range("M116").Value = 721.7777999999998
$ws.Range("N116").Value = -11914.1816
$ws.Range("H137").Value = 4572.8
$ws.Range("I137").Value = 4590.4
$ws.Range("K137").Value = 13771.2
$ws.Range("M137").Value = -11221.2
$ws.Range("H138").Value = 3518.2056
$ws.Range("I138").Value = 1580.8125
$ws.Range("J138").Value = 7238
$ws.Range("K138").Value = 4742.4375
$ws.Range("L138").Value = 21714
$ws.Range("M138").Value = 397.5625
$ws.Range("N138").Value = -31994
$ws.Range("H141").Value = 562957
$ws.Range("I141").Value = 1563.2354
$ws.Range("K141").Value = 4689.706200000001
$ws.Range("M141").Value = 490.2937999999995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5557234.5
$ws.Range("I2").Value = 12500822
$ws.Range("J2").Value = 2364.4
$ws.Range("K2").Value = 12500822
$ws.Range("L2").Value = 2364.4
$ws.Range("M2").Value = -12500709
$ws.Range("N2").Value = -2590.4
$ws.Range("H32").Value = 3818.8293
$ws.Range("I32").Value = 2558.2056
$ws.Range("K32").Value = 2558.2056
$ws.Range("M32").Value = -2271.2056
$ws.Range("H45").Value = 1516.1082
$ws.Range("I45").Value = 978.69696
$ws.Range("K45").Value = 978.69696
$ws.Range("M45").Value = -601.69696
$ws.Range("H116").Value = 5557234.5
$ws.Range("I116").Value = 12500822
$ws.Range("J116").Value = 2364.4
$ws.Range("K116").Value = 12500822
$ws.Range("L116").Value = 2364.4
$ws.Range("M116").Value = -12498528
$ws.Range("N116").Value = -6952.4
$ws.Range("H122").Value = 1855.8125
$ws.Range("I122").Value = 1201.3
$ws.Range("K122").Value = 3603.9
$ws.Range("M122").Value = -1153.9
$ws.Range("H135").Value = 34581.5
$ws.Range("J135").Value = 34581.5
$ws.Range("L135").Value = 34581.5
$ws.Range("N135").Value = -44721.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5557234.5
$ws.Range("I3").Value = 12500822
$ws.Range("J3").Value = 2364.4
$ws.Range("K3").Value = 12500822
$ws.Range("L3").Value = 2364.4
$ws.Range("M3").Value = -12500708
$ws.Range("N3").Value = -2592.4
$ws.Range("H134").Value = 3115.6052
$ws.Range("I134").Value = 2773.9678
$ws.Range("K134").Value = 8321.903399999999
$ws.Range("M134").Value = -5786.903399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2116.408
$ws.Range("I31").Value = 1340.5227
$ws.Range("J31").Value = 3183.25
$ws.Range("K31").Value = 1340.5227
$ws.Range("L31").Value = 3183.25
$ws.Range("M31").Value = -1045.5227
$ws.Range("N31").Value = -3773.25
$ws.Range("H34").Value = 2116.408
$ws.Range("I34").Value = 1340.5227
$ws.Range("J34").Value = 3183.25
$ws.Range("K34").Value = 1340.5227
$ws.Range("L34").Value = 3183.25
$ws.Range("M34").Value = -1138.5227
$ws.Range("N34").Value = -3587.25
$ws.Range("H94").Value = 3800.1667
$ws.Range("I94").Value = 5514.8887
$ws.Range("J94").Value = 2085.4443
$ws.Range("K94").Value = 5514.8887
$ws.Range("L94").Value = 2085.4443
$ws.Range("M94").Value = -5063.8887
$ws.Range("N94").Value = -2987.4443
$ws.Range("H105").Value = 3563
$ws.Range("I105").Value = 3118.9092
$ws.Range("K105").Value = 3118.9092
$ws.Range("M105").Value = -1371.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 202.4375
$ws.Range("I12").Value = 23
$ws.Range("K12").Value = 69
$ws.Range("M12").Value = 104
$ws.Range("H17").Value = 776
$ws.Range("I17").Value = 380
$ws.Range("J17").Value = 875
$ws.Range("K17").Value = 1140
$ws.Range("L17").Value = 2625
$ws.Range("M17").Value = -971
$ws.Range("N17").Value = -2963
$ws.Range("H68").Value = 2221.7917
$ws.Range("I68").Value = 743.9643
$ws.Range("J68").Value = 3162.2273
$ws.Range("K68").Value = 2231.8929
$ws.Range("L68").Value = 9486.6819
$ws.Range("M68").Value = -1420.8929
$ws.Range("N68").Value = -11108.6819
$ws.Range("H71").Value = 2221.7917
$ws.Range("I71").Value = 743.9643
$ws.Range("J71").Value = 3162.2273
$ws.Range("K71").Value = 6695.678699999999
$ws.Range("L71").Value = 28460.0457
$ws.Range("M71").Value = -2639.678699999999
$ws.Range("N71").Value = -36572.0457
$ws.Range("H112").Value = 2876.6667
$ws.Range("I112").Value = 725
$ws.Range("J112").Value = 3491.4285
$ws.Range("K112").Value = 2175
$ws.Range("L112").Value = 10474.2855
$ws.Range("M112").Value = -1067
$ws.Range("N112").Value = -12690.2855
$ws.Range("H129").Value = 29911.777
$ws.Range("I129").Value = 4226.6665
$ws.Range("J129").Value = 42754.332
$ws.Range("K129").Value = 12679.9995
$ws.Range("L129").Value = 128262.996
$ws.Range("M129").Value = -7679.999500000002
$ws.Range("N129").Value = -138262.996
$ws.Range("H131").Value = 1657.3158
$ws.Range("J131").Value = 1297.3043
$ws.Range("L131").Value = 3891.9129
$ws.Range("N131").Value = -13971.9129
$ws.Range("H132").Value = 1385.8
$ws.Range("J132").Value = 2215.8333
$ws.Range("L132").Value = 19942.4997
$ws.Range("N132").Value = -25002.4997
$ws.Range("H134").Value = 2159.4
$ws.Range("I134").Value = 1020.8
$ws.Range("J134").Value = 3298
$ws.Range("K134").Value = 3062.4
$ws.Range("L134").Value = 9894
$ws.Range("M134").Value = 2007.6
$ws.Range("N134").Value = -20034
$ws.Range("H139").Value = 5818588.5
$ws.Range("J139").Value = 13805.3
$ws.Range("L139").Value = 41415.89999999999
$ws.Range("N139").Value = -51695.89999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2731.8
$ws.Range("I126").Value = 1453
$ws.Range("J126").Value = 4650
$ws.Range("K126").Value = 4359
$ws.Range("L126").Value = 13950
$ws.Range("M126").Value = -1889
$ws.Range("N126").Value = -18890

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1615.25
$ws.Range("I100").Value = 799.44446
$ws.Range("J100").Value = 2282.7273
$ws.Range("K100").Value = 799.44446
$ws.Range("L100").Value = 2282.7273
$ws.Range("M100").Value = -258.44446
$ws.Range("N100").Value = -3364.7273
$ws.Range("H132").Value = 1701
$ws.Range("I132").Value = 1110.1923
$ws.Range("K132").Value = 3330.5769
$ws.Range("M132").Value = -800.5769
$ws.Range("H136").Value = 2962.4062
$ws.Range("I136").Value = 2581.4546
$ws.Range("K136").Value = 7744.3638
$ws.Range("M136").Value = -5194.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 897.2143
$ws.Range("I81").Value = 644.875
$ws.Range("J81").Value = 1233.6666
$ws.Range("K81").Value = 1289.75
$ws.Range("L81").Value = 2467.3332
$ws.Range("M81").Value = -228.75
$ws.Range("N81").Value = -4589.3332
$ws.Range("H84").Value = 897.2143
$ws.Range("I84").Value = 644.875
$ws.Range("J84").Value = 1233.6666
$ws.Range("K84").Value = 6448.75
$ws.Range("L84").Value = 12336.666
$ws.Range("M84").Value = -1144.75
$ws.Range("N84").Value = -22944.666
$ws.Range("H94").Value = 20297.5
$ws.Range("J94").Value = 20297.5
$ws.Range("L94").Value = 20297.5
$ws.Range("M94").Value = -22099.5
$ws.Range("H122").Value = 324569.4
$ws.Range("I122").Value = 501262.6
$ws.Range("J122").Value = 3309.0908
$ws.Range("K122").Value = 1503787.8
$ws.Range("L122").Value = 9927.2724
$ws.Range("M122").Value = -1501337.8
$ws.Range("N122").Value = -14827.2724
$ws.Range("H126").Value = 4547491.5
$ws.Range("I126").Value = 1677.2307
$ws.Range("J126").Value = 11113668
$ws.Range("K126").Value = 5031.6921
$ws.Range("L126").Value = 33341004
$ws.Range("M126").Value = -2561.6921
$ws.Range("N126").Value = -33345944
$ws.Range("H136").Value = 1528.1316
$ws.Range("I136").Value = 679
$ws.Range("J136").Value = 4712.375
$ws.Range("K136").Value = 2037
$ws.Range("L136").Value = 14137.125
$ws.Range("M136").Value = 513
$ws.Range("N136").Value = -19237.125
